$d = $word.ActiveDocument

# Each of these paragraphs currently holds the id value split across three
# runs: a "<id>" run, a plain id-text run, and a "</id>" run. The edit
# collapses each trio into a single run (keeping the "<id>"/"</id>" run's
# formatting) whose text is the full "<id>...</id>" string.
$ids = @("p059v_6", "p060r_1", "p060r_2", "p060r_3")

foreach ($id in $ids) {
    $openTag = "<id>"
    $closeTag = "</id>"
    $fullText = $openTag + $id + $closeTag

    # Locate the whole "<id>xxx</id>" span so we know where it starts/ends.
    $rngFull = $d.Content
    $rngFull.Find.ClearFormatting()
    $foundFull = $rngFull.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($foundFull) {
        $tagStart = $rngFull.Start
        $tagOpenEnd = $tagStart + $openTag.Length

        # Remove the id-text run and the closing "</id>" run (everything
        # after the opening "<id>" run), then append their combined text
        # back onto the still-surviving opening run so it becomes a single
        # run carrying the "<id>...</id>" text.
        $delRange = $d.Range($tagOpenEnd, $rngFull.End)
        $delRange.Text = ""

        $insertRange = $d.Range($tagOpenEnd, $tagOpenEnd)
        $insertRange.InsertAfter($id + $closeTag)
    }
}
